$d = $word.ActiveDocument

# The paragraph contains a field (begin/instrText " m:endcommentblock "/end)
# followed by a bold red run of text. The commit replaces the field with a
# plain literal-text run "{m:endcommentblock}" and prefixes the message run
# with "    <---".

$ooxml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p w:rsidP="00F5495F" w:rsidR="0017752F" w:rsidRDefault="0017752F">' +
'<w:r><w:t xml:space="preserve">{m:endcommentblock}</w:t></w:r>' +
'<w:r><w:rPr><w:b w:val="on"/><w:color w:val="FF0000"/></w:rPr>' +
'<w:t xml:space="preserve">    &lt;---Invalid block: Unexpected tag m:endcommentblock at this location</w:t></w:r>' +
'</w:p>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData>' +
'</pkg:part>' +
'</pkg:package>'

$d.Content.InsertXML($ooxml)
